# [EGSVC-42] - UI : Working on the Creating the agreement feature
#
# eisTestData.xlsx edit:
#   - assignmentDetails!G2 ("Position" for assignment1, ENGINEERING /
#     Assistant Engineer row): bump the position code suffix from
#     ENG_Assistant Engineer_1 -> ENG_Assistant Engineer_5
#   - the active sheet/selection moves from employeeDetails!G1 to
#     assignmentDetails!G2 (i.e. assignmentDetails becomes the tab-selected
#     sheet, with employeeDetails' own cursor left parked on F2)

$wb = $excel.ActiveWorkbook

$employeeDetails    = $wb.Worksheets.Item("employeeDetails")
$assignmentDetails  = $wb.Worksheets.Item("assignmentDetails")

# Update the assignment "Position" code used for the Assistant Engineer row.
$assignmentDetails.Range("G2").Value = "ENG_Assistant Engineer_5"

# Leave the employeeDetails cursor on F2 ...
$employeeDetails.Range("F2").Select() | Out-Null

# ... and make assignmentDetails (with the cell we just edited selected)
# the active sheet/tab, cursor on G2.
$assignmentDetails.Activate()
$assignmentDetails.Range("G2").Select() | Out-Null
